$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5: "Escrita científica" / "Writing Reproducible Research Papers with R
# Markdown" / link to resulumit.com workshop slides
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Escrita científica"
$ws.Range("C5").Value = "Writing Reproducible Research Papers with R Markdown"

$ws.Hyperlinks.Add($ws.Range("D5"), "https://resulumit.com/teaching/rmd_workshop.html#1", "", "", "https://resulumit.com/teaching/rmd_workshop.html#1")
$ws.Range("D5").Font.Name = "Calibri"
$ws.Range("D5").Font.Size = 11
$ws.Range("D5").Font.Underline = 2
$ws.Range("D5").Font.Color = 255

# ---------------------------------------------------------------------------
# Row 6: "RMarkdown" / "R Markdown Cheatsheet" / link to rstudio cheatsheet
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "RMarkdown"
$ws.Range("C6").Value = "R Markdown Cheatsheet"

$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.rstudio.com/wp-content/uploads/2015/02/rmarkdown-cheatsheet.pdf", "", "", "https://www.rstudio.com/wp-content/uploads/2015/02/rmarkdown-cheatsheet.pdf")
$ws.Range("D6").Font.Name = "Calibri"
$ws.Range("D6").Font.Size = 11
$ws.Range("D6").Font.Underline = 2
$ws.Range("D6").Font.Color = 255
